$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Test Sonucu"
$ws.Range("C2").Value = "PASSED"
$ws.Range("C3").Value = "FAILED"
